$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 3 (per diff: D3 2->3, F3 2->3, H3 36->46)
$ws.Range("D3").Value = 3
$ws.Range("F3").Value = 3
$ws.Range("H3").Value = 46

# Update the active selection cell from D5 to C3
$ws.Range("C3").Select() | Out-Null
